$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "mngr348764"
$ws.Range("B7").Value = "ydYhUvy"

$ws.Range("A7:B7").Font.Bold = $true
$ws.Range("A7:B7").Font.Size = 14
$ws.Range("A7:B7").Font.Name = "Arial"
$ws.Range("A7:B7").Font.Color = 0

$ws.Range("B7").Select()
